# Delete column AU. This shifts every cell from AU onward one column to the
# left (AV -> AU, AW -> AV, ... BK -> BJ), matching the target diff:
#   - AU1's old "MO" label is dropped; AV1..AZ1 ("UTMD","DE","DKNG","G","IT")
#     shift left into AU1..AY1 (formatting shifts with them).
#   - Row 2/3/4 values in columns AU..BK shift left by one column, and the
#     now-vacated last column (BK) disappears from the used range.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns("AU").Delete()
